$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.581.56'
$ws.Range("E2").Value = '  -0.17%  '
$ws.Range("D3").Value = '1.752.30'
$ws.Range("E3").Value = '  -0.43%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '324.43'
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("D7").Value = '0.4491'
$ws.Range("E7").Value = '  +2.95%  '
$ws.Range("E8").Value = '  -1.26%  '
$ws.Range("D9").Value = '0.07466'
$ws.Range("E9").Value = '  -1.34%  '
$ws.Range("D10").Value = '41.50'
$ws.Range("E10").Value = '  -1.41%  '
$ws.Range("E11").Value = '  -2.70%  '
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  +0.10%  '
$ws.Range("E13").Value = '  -0.58%  '
$ws.Range("D14").Value = '5.980'
$ws.Range("E14").Value = '  -1.60%  '
$ws.Range("D15").Value = '7.148'
$ws.Range("E15").Value = '  -1.19%  '
$ws.Range("D16").Value = '1.756.06'
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").Value = '93.63'
$ws.Range("E17").Value = '  +0.99%  '
$ws.Range("D18").Value = '0.00001056'
$ws.Range("E18").Value = '  -1.18%  '
$ws.Range("D19").Value = '0.06381'
$ws.Range("E19").Value = '  -0.52%  '
$ws.Range("D21").Value = '17.15'
$ws.Range("E21").Value = '  +0.27%  '
$ws.Range("D22").Value = '5.732'
$ws.Range("E22").Value = '  -2.05%  '
$ws.Range("D23").Value = '27.634.20'
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("E24").Value = '  -0.67%  '
$ws.Range("D25").Value = '2.088'
$ws.Range("E25").Value = '  -0.46%  '
$ws.Range("D26").Value = '165.88'
$ws.Range("E26").Value = '  +1.87%  '
$ws.Range("D27").Value = '20.14'
$ws.Range("E27").Value = '  -1.64%  '
$ws.Range("D28").Value = '1.957.24'
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("D29").Value = '2.091'
$ws.Range("E29").Value = '  -2.56%  '
$ws.Range("D30").Value = '125.61'
$ws.Range("E30").Value = '  -0.37%  '
$ws.Range("D31").Value = '1.093'
$ws.Range("E31").Value = '  -0.69%  '
$ws.Range("D32").Value = '0.09177'
$ws.Range("E32").Value = '  +2.08%  '
$ws.Range("D33").Value = '3.663'
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("E34").Value = '  -1.46%  '
$ws.Range("D35").Value = '0.02286'
$ws.Range("E35").Value = '  -0.88%  '
$ws.Range("E36").Value = '  -3.89%  '
$ws.Range("D37").Value = '0.2093'
$ws.Range("E37").Value = '  -0.76%  '
$ws.Range("D38").Value = '0.06023'
$ws.Range("E38").Value = '  +0.11%  '
$ws.Range("D39").Value = '0.6287'
$ws.Range("E39").Value = '  -2.35%  '
$ws.Range("D40").Value = '4.923'
$ws.Range("D41").Value = '1.182'
$ws.Range("E41").Value = '  -0.37%  '
$ws.Range("D42").Value = '1.389'
$ws.Range("E42").Value = '  -0.88%  '
$ws.Range("D43").Value = '7.760'
$ws.Range("E43").Value = '  -1.74%  '
$ws.Range("D44").Value = '13.09'
$ws.Range("E44").Value = '  -1.73%  '
$ws.Range("D45").Value = '3.715'
$ws.Range("E45").Value = '  -0.08%  '
$ws.Range("D46").Value = '0.5866'
$ws.Range("E46").Value = '  -1.24%  '
$ws.Range("D47").Value = '122.26'
$ws.Range("E47").Value = '  -0.03%  '
$ws.Range("D48").Value = '1.935'
$ws.Range("E48").Value = '  -2.73%  '
$ws.Range("D49").Value = '0.06887'
$ws.Range("E49").Value = '  +0.09%  '
$ws.Range("D50").Value = '1.129'
$ws.Range("E50").Value = '  -3.60%  '
$ws.Range("D51").Value = '71.59'
$ws.Range("E51").Value = '  -1.94%  '
